# Insert a new data row at row 230 (pushes existing rows 230..334 down to 231..335)
# and populate it with the new record. This mirrors the author's edit which added
# one additional weekly "Poroto verde" observation in the middle of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(230).Insert()

$ws.Cells.Item(230, 1).Value  = 3
$ws.Cells.Item(230, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(230, 3).Value  = "Coquimbo"
$ws.Cells.Item(230, 4).Value  = 44636
$ws.Cells.Item(230, 5).Value  = 5
$ws.Cells.Item(230, 6).Value  = 100112031
$ws.Cells.Item(230, 7).Value  = "Poroto verde"
$ws.Cells.Item(230, 8).Value  = "Magnum"
$ws.Cells.Item(230, 9).Value  = "Primera"
$ws.Cells.Item(230, 10).Value = 103
$ws.Cells.Item(230, 11).Value = 23000
$ws.Cells.Item(230, 12).Value = 24000
$ws.Cells.Item(230, 13).Value = 23340
$ws.Cells.Item(230, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(230, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(230, 16).Value = 934
$ws.Cells.Item(230, 17).Value = 25
$ws.Cells.Item(230, 18).Value = "Hortaliza"

$ws.Cells.Item(230, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
